# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Adds a new "2506" period row to the worker's account-statement table
# (inserted between the existing "2505" and "2506" rows, which become the
# "2507" and "2505" rows respectively), bumps the period count and the
# total overdue amount accordingly, and reorders the "Novedad de
# Ingreso/Retiro" header columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Totals: one more period (3 instead of 2) at 56940 each ---
$ws.Range("E11").Value = 170820
$ws.Range("F13").Value = 3

# --- Header row: swap "Novedad de Ingreso" / "Novedad de Retiro" order ---
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# --- Insert a new detail row right after the first data row (row 17), ---
# --- copying row 16's formatting so the new row matches the existing ---
# --- table styling exactly.                                          ---
$ws.Rows.Item(17).Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 16 keeps the worker but now reports period 2507
$ws.Range("E16").Value = "2507"

# New row 17: same worker, period 2506
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1050964669"
$ws.Range("D17").Value = "ALDO ENRIQUE GONZALEZ DIAZ"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Former row 17 (now row 18, pushed down by the insert) reports period 2505
$ws.Range("E18").Value = "2505"
